$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Achievements")

# My Preciouses (row 5): Estimated Level Acquired 499 -> 567, Items Needed 2800 -> 3000
$ws.Range("C5").Value = 567
$ws.Range("D5").Value = 3000

# Exterminator (row 6): Estimated Level Acquired 140 -> 215, Items Needed 50 -> 100
$ws.Range("C6").Value = 215
$ws.Range("D6").Value = 100

# Dragon Slayer (row 7): Estimated Level Acquired 215 -> 350, Items Needed 100 -> 200
$ws.Range("C7").Value = 350
$ws.Range("D7").Value = 200

# Beast Master (row 8): Estimated Level Acquired 425 -> 790, Items Needed 300 -> 550
$ws.Range("C8").Value = 790
$ws.Range("D8").Value = 550

# Boulder Breaker (row 10): Estimated Level Acquired 200 -> 250, Items Needed 85 -> 125
$ws.Range("C10").Value = 250
$ws.Range("D10").Value = 125

# Rock n Roller (row 11): Estimated Level Acquired 400 -> 675, Items Needed 230 -> 450
$ws.Range("C11").Value = 675
$ws.Range("D11").Value = 450

# Hoarder (row 14): Estimated Level Acquired text "350?" -> "400?"
$ws.Range("C14").Value = "400?"

# PUZL Master (row 30): Estimated Level Acquired 500 -> 790
$ws.Range("C30").Value = 790

# Restore the active selection to D5, matching the saved sheet view state
$ws.Range("D5").Select()
